$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, pushing the existing row 55 (and below) down to row 56.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new data record.
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 44939
$ws.Range("D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 100112030
$ws.Range("G55").Value = "Poroto granado"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 70
$ws.Range("K55").Value = 45000
$ws.Range("L55").Value = 45000
$ws.Range("M55").Value = 45000
$ws.Range("N55").Value = "$/saco 25 kilos"
$ws.Range("O55").Value = "Región Metropolitana"
$ws.Range("P55").Value = 1800
$ws.Range("Q55").Value = 25
$ws.Range("R55").Value = "Hortaliza"
